# Refresh market-price-derived columns (H:N) on each profession sheet.
# Values come from the scheduled market-data runner; table formulas are not
# used in this workbook so cells are written directly.
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 20972.445
$ws.Range("I28").Value = 5615.8335
$ws.Range("J28").Value = 51685.668
$ws.Range("K28").Value = 5615.8335
$ws.Range("L28").Value = 51685.668
$ws.Range("M28").Value = -5130.8335
$ws.Range("N28").Value = -52655.668
# Row 38
$ws.Range("H38").Value = 779.7619
$ws.Range("I38").Value = 189.33333
$ws.Range("J38").Value = 1567
$ws.Range("K38").Value = 567.99999
$ws.Range("L38").Value = 4701
$ws.Range("M38").Value = -195.99999
$ws.Range("N38").Value = -5445
# Row 61
$ws.Range("H61").Value = 595.5294
$ws.Range("I61").Value = 394.93332
$ws.Range("J61").Value = 2100
$ws.Range("K61").Value = 1184.79996
$ws.Range("L61").Value = 6300
$ws.Range("M61").Value = -1012.79996
$ws.Range("N61").Value = -6644
# Row 111
$ws.Range("H111").Value = 8414.5
$ws.Range("I111").Value = 10029
$ws.Range("J111").Value = 6800
$ws.Range("K111").Value = 30087
$ws.Range("L111").Value = 20400
$ws.Range("M111").Value = -27020
$ws.Range("N111").Value = -26534
# Row 112
$ws.Range("H112").Value = 13376.429
$ws.Range("J112").Value = 13376.429
$ws.Range("L112").Value = 40129.287
$ws.Range("N112").Value = -42345.287
# Row 118
$ws.Range("H118").Value = 2908
$ws.Range("I118").Value = 2582.8572
$ws.Range("J118").Value = 3666.6667
$ws.Range("K118").Value = 7748.571599999999
$ws.Range("L118").Value = 11000.0001
$ws.Range("M118").Value = -6091.571599999999
$ws.Range("N118").Value = -14314.0001
# Row 132
$ws.Range("H132").Value = 4922.1
$ws.Range("I132").Value = 5762.923
$ws.Range("J132").Value = 3360.5715
$ws.Range("K132").Value = 17288.769
$ws.Range("L132").Value = 10081.7145
$ws.Range("M132").Value = -14758.769
$ws.Range("N132").Value = -15141.7145
# Row 137
$ws.Range("H137").Value = 1546.4445
$ws.Range("I137").Value = 1399.6
$ws.Range("J137").Value = 1730
$ws.Range("K137").Value = 4198.799999999999
$ws.Range("L137").Value = 5190
$ws.Range("M137").Value = -1648.799999999999
$ws.Range("N137").Value = -10290

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 200
$ws.Range("J4").Value = 200
$ws.Range("L4").Value = 200
$ws.Range("N4").Value = -432
# Row 6
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 500
$ws.Range("K6").Value = 500
$ws.Range("M6").Value = -327
# Row 32
$ws.Range("H32").Value = 1555082.5
$ws.Range("I32").Value = 2029601.2
$ws.Range("J32").Value = 19874.412
$ws.Range("K32").Value = 2029601.2
$ws.Range("L32").Value = 19874.412
$ws.Range("M32").Value = -2029314.2
$ws.Range("N32").Value = -20448.412
# Row 45
$ws.Range("H45").Value = 2253.0557
$ws.Range("I45").Value = 1885.091
$ws.Range("J45").Value = 2831.2856
$ws.Range("K45").Value = 1885.091
$ws.Range("L45").Value = 2831.2856
$ws.Range("M45").Value = -1508.091
$ws.Range("N45").Value = -3585.2856
# Row 122
$ws.Range("H122").Value = 146400
$ws.Range("I122").Value = 1000000
$ws.Range("K122").Value = 3000000
$ws.Range("M122").Value = -2997550

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 466.66666
$ws.Range("I22").Value = 468.75
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 468.75
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -295.75
$ws.Range("N22").Value = -796

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 67.333336
$ws.Range("I7").Value = 50.857143
$ws.Range("J7").Value = 125
$ws.Range("K7").Value = 50.857143
$ws.Range("L7").Value = 125
$ws.Range("M7").Value = 62.142857
$ws.Range("N7").Value = -351
# Row 22
$ws.Range("H22").Value = 658.3333
$ws.Range("I22").Value = 562.5
$ws.Range("J22").Value = 850
$ws.Range("K22").Value = 562.5
$ws.Range("L22").Value = 850
$ws.Range("M22").Value = -212.5
$ws.Range("N22").Value = -1550
# Row 31
$ws.Range("H31").Value = 5136.1772
$ws.Range("I31").Value = 1383.25
$ws.Range("J31").Value = 6923.2856
$ws.Range("K31").Value = 1383.25
$ws.Range("L31").Value = 6923.2856
$ws.Range("M31").Value = -1088.25
$ws.Range("N31").Value = -7513.2856
# Row 34
$ws.Range("H34").Value = 5136.1772
$ws.Range("I34").Value = 1383.25
$ws.Range("J34").Value = 6923.2856
$ws.Range("K34").Value = 1383.25
$ws.Range("L34").Value = 6923.2856
$ws.Range("M34").Value = -1181.25
$ws.Range("N34").Value = -7327.2856

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 621.4286
$ws.Range("I86").Value = 620
$ws.Range("J86").Value = 622.5
$ws.Range("K86").Value = 1860
$ws.Range("L86").Value = 1867.5
$ws.Range("M86").Value = -674
$ws.Range("N86").Value = -4239.5
# Row 89
$ws.Range("H89").Value = 621.4286
$ws.Range("I89").Value = 620
$ws.Range("J89").Value = 622.5
$ws.Range("K89").Value = 5580
$ws.Range("L89").Value = 5602.5
$ws.Range("M89").Value = 348
$ws.Range("N89").Value = -17458.5
# Row 131
$ws.Range("H131").Value = 4129.75
$ws.Range("J131").Value = 6240.0435
$ws.Range("L131").Value = 18720.1305
$ws.Range("N131").Value = -28800.1305
# Row 137
$ws.Range("H137").Value = 30200.55
$ws.Range("I137").Value = 7170.6313
$ws.Range("J137").Value = 51037.145
$ws.Range("K137").Value = 21511.8939
$ws.Range("L137").Value = 153111.435
$ws.Range("M137").Value = -16411.8939
$ws.Range("N137").Value = -163311.435
# Row 140
$ws.Range("H140").Value = 2550.375
$ws.Range("I140").Value = 2267.1428
$ws.Range("K140").Value = 6801.428400000001
$ws.Range("M140").Value = -1621.428400000001

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 35900
$ws.Range("I122").Value = 52500
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 157500
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -155050
$ws.Range("N122").Value = -13000

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 20510.2
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 25362.75
$ws.Range("K22").Value = 1100
$ws.Range("L22").Value = 25362.75
$ws.Range("M22").Value = -805
$ws.Range("N22").Value = -25952.75
# Row 27
$ws.Range("H27").Value = 20510.2
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 25362.75
$ws.Range("K27").Value = 1100
$ws.Range("L27").Value = 25362.75
$ws.Range("M27").Value = -993
$ws.Range("N27").Value = -25576.75
# Row 40
$ws.Range("H40").Value = 3500
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2864
# Row 132
$ws.Range("H132").Value = 3196.6775
$ws.Range("J132").Value = 4356.7144
$ws.Range("L132").Value = 13070.1432
$ws.Range("N132").Value = -18130.1432

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1186.5
$ws.Range("I126").Value = 1238.1666
$ws.Range("J126").Value = 1083.1666
$ws.Range("K126").Value = 3714.4998
$ws.Range("L126").Value = 3249.4998
$ws.Range("M126").Value = -1244.4998
$ws.Range("N126").Value = -8189.4998
# Row 132
$ws.Range("H132").Value = 6631386
$ws.Range("I132").Value = 2365.2903
$ws.Range("K132").Value = 7095.8709
$ws.Range("M132").Value = -4565.8709
